# Appends newly-logged sensor readings (2026-02-06) to the PIR, Humidity,
# and Temperature sheets, continuing the existing Bathroom sensor log.
$wb = $excel.ActiveWorkbook

function Append-SensorRows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        $Rows
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $endRow = $StartRow + $Rows.Count - 1
    $ws.Range("A" + $StartRow + ":F" + $endRow).NumberFormat = "@"
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $values = $Rows[$i]
        for ($c = 0; $c -lt $values.Count; $c++) {
            $ws.Cells.Item($r, $c + 1).Value = $values[$c]
        }
    }
}

# --- PIR sheet: new rows 353-365 ---
$pirRows = @(
    @('2026-02-06','10:09:17','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:09:19','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:09:23','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:09:29','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:09:34','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:09:39','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:09:44','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:09:49','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:09:54','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:09:59','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:10:04','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:10:08','10:00','Bathroom','No Motion','Inactive'),
    @('2026-02-06','10:10:14','10:00','Bathroom','No Motion','Inactive')
)
Append-SensorRows "PIR" 353 $pirRows

# --- Humidity sheet: new rows 236-246 ---
$humidityRows = @(
    @('2026-02-06','10:09:18','10:00','Bathroom','68.9%','Active'),
    @('2026-02-06','10:09:22','10:00','Bathroom','68.9%','Active'),
    @('2026-02-06','10:09:27','10:00','Bathroom','68.9%','Active'),
    @('2026-02-06','10:09:32','10:00','Bathroom','68.9%','Active'),
    @('2026-02-06','10:09:37','10:00','Bathroom','68.8%','Active'),
    @('2026-02-06','10:09:42','10:00','Bathroom','68.9%','Active'),
    @('2026-02-06','10:09:47','10:00','Bathroom','68.9%','Active'),
    @('2026-02-06','10:09:52','10:00','Bathroom','68.8%','Active'),
    @('2026-02-06','10:09:57','10:00','Bathroom','68.7%','Active'),
    @('2026-02-06','10:10:02','10:00','Bathroom','68.6%','Active'),
    @('2026-02-06','10:10:12','10:00','Bathroom','68.7%','Active')
)
Append-SensorRows "Humidity" 236 $humidityRows

# --- Temperature sheet: new rows 236-246 ---
$temperatureRows = @(
    @('2026-02-06','10:09:19','10:00','Bathroom','28.1C','Active'),
    @('2026-02-06','10:09:23','10:00','Bathroom','28.1C','Active'),
    @('2026-02-06','10:09:28','10:00','Bathroom','28.1C','Active'),
    @('2026-02-06','10:09:33','10:00','Bathroom','28.1C','Active'),
    @('2026-02-06','10:09:38','10:00','Bathroom','28.1C','Active'),
    @('2026-02-06','10:09:43','10:00','Bathroom','28.1C','Active'),
    @('2026-02-06','10:09:48','10:00','Bathroom','28.1C','Active'),
    @('2026-02-06','10:09:53','10:00','Bathroom','28.1C','Active'),
    @('2026-02-06','10:09:58','10:00','Bathroom','28.1C','Active'),
    @('2026-02-06','10:10:03','10:00','Bathroom','28.0C','Active'),
    @('2026-02-06','10:10:13','10:00','Bathroom','28.0C','Active')
)
Append-SensorRows "Temperature" 236 $temperatureRows

